$d = $word.ActiveDocument

# "main: Update docs Year" - bump the year in the "Nomor : 470/{noReg}/14.2013/2022"
# line to 2023.
$d.Content.Find.Execute("2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, "2023", 2)
